$d = $word.ActiveDocument

# --- Step 1: flip the 8 "Failure" (red) results to "True" (green) ---
# Row indices below are 1-based Table.Rows indices (row 1 is the header row).
$rowsToChange = @(2, 3, 5, 6, 8, 11, 12, 13)
foreach ($r in $rowsToChange) {
    $t = $d.Tables.Item(1)
    $cell = $t.Cell($r, 2)
    $word1 = $cell.Range.Words.Item(1)
    $word1.Font.Color = 5287936   # RGB(0x00, 0xB0, 0x50) -> w:color 00B050
    $word1.Text = "True"
}

# --- Step 2: append three new rows to the table ---
$newRowsData = @(
    @("All data is stored in the Model part of the application", "True"),
    @("All visuals are stored in the View part of the application", "True"),
    @("All logic is stored in the Controller part of the application", "True")
)

foreach ($pair in $newRowsData) {
    $t = $d.Tables.Item(1)
    $newRow = $t.Rows.Add()
    $cell1 = $newRow.Cells.Item(1)
    $cell2 = $newRow.Cells.Item(2)
    $cell1.Range.Text = $pair[0]
    $cell2.Range.Text = $pair[1]
    $cell2.Range.Font.Color = 5287936
}

Write-Host "done"
